# Applies commit "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# - Updates "Periodo Mora" from 2508 to 2509 for existing workers
# - Inserts a new worker row (ANGIE PAOLA MORENO ARIAS) in the debtor table
# - Updates the aggregate "Cant. Trabajadores" and "VALOR MORA" totals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row right below the first worker (STELLA ANGULO CARMONA, row 16),
#    pushing the existing second worker (TIVEL ESTEPHANI BATISTA MARTINEZ) and the
#    signature block further down.
$ws.Rows.Item(17).Insert()

# Copy formatting from the row above so the new row matches the table's borders/fills.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 2. Update "Periodo Mora" for the pre-existing rows (2508 -> 2509).
$ws.Range("E16").Value = "2509"
$ws.Range("E18").Value = "2509"

# 3. Fill in the new worker's row.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1007263396"
$ws.Range("D17").Value = "ANGIE PAOLA MORENO ARIAS"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# 4. Update the summary fields.
$ws.Range("C13").Value = 3
$ws.Range("E11").Value = 316940
